# "Fix typo in test file"
#
# testUnit.xlsx contained a misspelled sheet name ("TabelSheet") and a
# misspelled built-in cell-style name ("Normal" should read "Standard" in
# the German-localised original). This script corrects the sheet name
# (the substantive fix) and also nudges the per-sheet selections back to
# the state captured by the re-saved workbook: the active tab moves from
# "DataTypeSheet" to the corrected "TableSheet", and "NumericSheet"'s
# remembered selection moves from B38 to C7.

$wb = $excel.ActiveWorkbook

# --- the actual typo fix: TabelSheet -> TableSheet -------------------------
$wsTable = $wb.Sheets.Item("TabelSheet")
$wsTable.Name = "TableSheet"

# --- related typo: built-in style "Normal" -> "Standard" -------------------
try {
    $wb.Styles.Item("Normal").Name = "Standard"
} catch {
    # older/headless hosts may not expose a writable Styles collection;
    # ignore if unsupported.
}

# --- StringSheet: collapse the old A1:D13 block-selection down to A1 -------
$ws1 = $wb.Sheets.Item("StringSheet")
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null

# --- NumericSheet: selection moves from B38 to C7 ---------------------------
$ws2 = $wb.Sheets.Item("NumericSheet")
$ws2.Activate() | Out-Null
$ws2.Range("C7").Select() | Out-Null

# --- TableSheet (renamed) keeps its own C6 selection, but becomes the ------
# --- active / last-selected tab (taking that status from DataTypeSheet) ----
$wsTable.Activate() | Out-Null
